$wb = $excel.ActiveWorkbook

# --- Sheet "data" (sheet1): add column AN (col 40) ---
$ws1 = $wb.Worksheets.Item("data")

# Header cell AN1: copy formatting from AM1 (style s="1"), then set text
$ws1.Range("AM1").Copy()
$ws1.Range("AN1").PasteSpecial(-4122)
$ws1.Range("AN1").Value = "25. 1. 2021"

# Data rows 2..67: numeric values for column AN
$data1 = @{
    2 = 0.08
    3 = 0.18
    4 = 0.74
    5 = 0.07000000000000001
    6 = 0.14
    7 = 0.79
    8 = 0.03
    9 = 0.13
    10 = 0.84
    11 = 0.08
    12 = 0.23
    13 = 0.6899999999999999
    14 = 0.2
    15 = 0.16
    16 = 0.64
    17 = 0.07000000000000001
    18 = 0.14
    19 = 0.79
    20 = 0.04
    21 = 0.13
    22 = 0.83
    23 = 0.18
    24 = 0.18
    25 = 0.64
    26 = 0.1
    27 = 0.22
    28 = 0.68
    29 = 0.14
    30 = 0.26
    31 = 0.6
    32 = 0.06
    33 = 0.14
    34 = 0.8
    35 = 0.01
    36 = 0.07000000000000001
    37 = 0.92
    38 = 0.18
    39 = 0.17
    40 = 0.65
    41 = 0.06
    42 = 0.18
    43 = 0.76
    44 = 0.6
    45 = 0.19
    46 = 0.21
    47 = 0.11
    48 = 0.55
    49 = 0.34
    50 = 0.02
    51 = 0.08
    52 = 0.9
    53 = 0.07000000000000001
    54 = 0.18
    55 = 0.75
    56 = 0
    57 = 0.11
    58 = 0.89
    59 = 0.11
    60 = 0.17
    61 = 0.72
    62 = 0.02
    63 = 0.08
    64 = 0.9
    65 = 0.06
    66 = 0.11
    67 = 0.83
}
foreach ($r in $data1.Keys) {
    $ws1.Cells.Item($r, 40).Value = $data1[$r]
}

# Footer row 68: update "aktualizace" date in the label
$ws1.Cells.Item(68, 1).Value = "Život během pandemie, Zasažení domácností, % respondentů celkově a ve skupinách, aktualizace 1. 2. 2022"

# --- Sheet "pocetR" (sheet2): add column AM (col 39) ---
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AM1: copy formatting from AL1 (style s="2"), then set text
$ws2.Range("AL1").Copy()
$ws2.Range("AM1").PasteSpecial(-4122)
$ws2.Range("AM1").Value = "25. 1. 2021"

# Data rows 2..23: numeric values for column AM
$data2 = @{
    2 = 1560
    3 = 740
    4 = 124
    5 = 483
    6 = 213
    7 = 707
    8 = 114
    9 = 100
    10 = 639
    11 = 729
    12 = 529
    13 = 302
    14 = 429
    15 = 1131
    16 = 141
    17 = 308
    18 = 1111
    19 = 272
    20 = 83
    21 = 254
    22 = 140
    23 = 82
}
foreach ($r in $data2.Keys) {
    $ws2.Cells.Item($r, 39).Value = $data2[$r]
}

# Footer row 24: empty trailing placeholder cell AM24 (matches B24:AL24 pattern)
$ws2.Cells.Item(24, 39).Style = $ws2.Cells.Item(24, 38).Style

# Footer row 24: update "aktualizace" date in the label
$ws2.Cells.Item(24, 1).Value = "Život během pandemie, Zasažení domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 2. 2022"

